# Update invoice currency and payment-terms text.
#
#   Currency:      "USD"             -> "EUR"
#   Payment Terms: "due in 45 days"  -> "45 days"
#
# Both values live in row 2, underneath the "Currency" / "Payment Terms"
# headers in row 1, on every worksheet that has them ("Simple Fields" and
# "Simple Fields - Formatted"); the "Items" sheets don't carry these
# columns, so they're naturally skipped by the header lookup below.

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    $usedRange = $ws.UsedRange
    $headerRow = $usedRange.Rows.Item(1)
    $lastCol = $headerRow.Columns.Count

    for ($col = 1; $col -le $lastCol; $col++) {
        $header = $ws.Cells.Item(1, $col).Value()
        $dataCell = $ws.Cells.Item(2, $col)

        if ($header -eq "Currency") {
            if ($dataCell.Value() -eq "USD") {
                $dataCell.Value = "EUR"
            }
        }
        elseif ($header -eq "Payment Terms") {
            if ($dataCell.Value() -eq "due in 45 days") {
                $dataCell.Value = "45 days"
            }
        }
    }
}
